$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 4, 5 and 6 - their data is superseded by the
# new row 2/row 3 content set below.
$ws.Rows("4:6").Delete()

# Row 2: new "044/FES VILLE" / "Direction régionale" entry (rappel avenant MA)
$ws.Range("A2").Value = "044/FES VILLE "
$ws.Range("B2").Value = "Direction régionale"
$ws.Range("C2").Value = "K5443645"
$ws.Range("D2").Value = "KHADIJA LALA"
$ws.Range("E2").Value = "non"
$ws.Range("F2").Value = "mensuelle"
$ws.Range("G2").Value = 15
$ws.Range("H2").Value = 10000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 10000
$ws.Range("O2").Value = 18500

# Row 3: blank identity columns, matching amounts to row 2
$ws.Range("A3").Value = " "
$ws.Range("B3").Value = " "
$ws.Range("C3").Value = " "
$ws.Range("D3").Value = " "
$ws.Range("E3").Value = " "
$ws.Range("F3").Value = " "
$ws.Range("G3").Value = " "
$ws.Range("H3").Value = 10000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 10000
$ws.Range("O3").Value = 18500
